# This script updates the TPM-derived NATMI ligand-receptor statistics
# for the Slit3 -> Robo2 interaction table (Sheet1) to reflect newly
# recomputed TPM values, per "update scripts wuth new tpm".
#
# Columns:
#   G/H  = Ligand average/total expression value (sending cluster)
#   I/J  = Ligand derived specificity (average/total)
#   K/L  = Receptor-expressing cells / detection rate (target cluster)
#   M/N  = Receptor average/total expression value (target cluster)
#   O/P  = Receptor derived specificity (average/total)
#   Q/R  = Edge average/total expression weight
#   S/T  = Edge average/total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.232451333333333
$ws.Range("H2").Value = 6.697354
$ws.Range("I2").Value = 0.01414074962829973
$ws.Range("J2").Value = 0.01414074962829973
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.892869333333334
$ws.Range("N2").Value = 8.678608
$ws.Range("O2").Value = 0.9698236995656885
$ws.Range("P2").Value = 0.9698236995656884
$ws.Range("Q2").Value = 6.458190000359111
$ws.Range("R2").Value = 58.123710003232
$ws.Range("S2").Value = 0.01371403411914978
$ws.Range("T2").Value = 0.01371403411914978
$ws.Range("G3").Value = 2.232451333333333
$ws.Range("H3").Value = 6.697354
$ws.Range("I3").Value = 0.01414074962829973
$ws.Range("J3").Value = 0.01414074962829973
$ws.Range("O3").Value = 0.006810193051573731
$ws.Range("P3").Value = 0.00681019305157373
$ws.Range("Q3").Value = 0.04535001638533333
$ws.Range("R3").Value = 0.408150147468
$ws.Range("S3").Value = 0.00009630123486269065
$ws.Range("T3").Value = 0.00009630123486269064
$ws.Range("G4").Value = 2.232451333333333
$ws.Range("H4").Value = 6.697354
$ws.Range("I4").Value = 0.01414074962829973
$ws.Range("J4").Value = 0.01414074962829973
$ws.Range("O4").Value = 0.02336610738273784
$ws.Range("P4").Value = 0.02336610738273783
$ws.Range("Q4").Value = 0.1555981371811111
$ws.Range("R4").Value = 1.40038323463
$ws.Range("S4").Value = 0.0003304142742872617
$ws.Range("T4").Value = 0.0003304142742872616
$ws.Range("I5").Value = 0.8099327614075106
$ws.Range("J5").Value = 0.8099327614075106
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.892869333333334
$ws.Range("N5").Value = 8.678608
$ws.Range("O5").Value = 0.9698236995656885
$ws.Range("P5").Value = 0.9698236995656884
$ws.Range("Q5").Value = 369.902572224112
$ws.Range("R5").Value = 3329.123150017008
$ws.Range("S5").Value = 0.785491987067686
$ws.Range("T5").Value = 0.7854919870676859
$ws.Range("I6").Value = 0.8099327614075106
$ws.Range("J6").Value = 0.8099327614075106
$ws.Range("O6").Value = 0.006810193051573731
$ws.Range("P6").Value = 0.00681019305157373
$ws.Range("S6").Value = 0.005515798463979353
$ws.Range("T6").Value = 0.005515798463979352
$ws.Range("I7").Value = 0.8099327614075106
$ws.Range("J7").Value = 0.8099327614075106
$ws.Range("O7").Value = 0.02336610738273784
$ws.Range("P7").Value = 0.02336610738273783
$ws.Range("S7").Value = 0.01892497587584527
$ws.Range("T7").Value = 0.01892497587584527
$ws.Range("G8").Value = 27.77415166666666
$ws.Range("H8").Value = 83.32245499999999
$ws.Range("I8").Value = 0.1759264889641896
$ws.Range("J8").Value = 0.1759264889641896
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.892869333333334
$ws.Range("N8").Value = 8.678608
$ws.Range("O8").Value = 0.9698236995656885
$ws.Range("P8").Value = 0.9698236995656884
$ws.Range("Q8").Value = 80.34699161584889
$ws.Range("R8").Value = 723.12292454264
$ws.Range("S8").Value = 0.1706176783788527
$ws.Range("T8").Value = 0.1706176783788526
$ws.Range("G9").Value = 27.77415166666666
$ws.Range("H9").Value = 83.32245499999999
$ws.Range("I9").Value = 0.1759264889641896
$ws.Range("J9").Value = 0.1759264889641896
$ws.Range("O9").Value = 0.006810193051573731
$ws.Range("P9").Value = 0.00681019305157373
$ws.Range("Q9").Value = 0.5642041169566667
$ws.Range("S9").Value = 0.001198093352731687
$ws.Range("T9").Value = 0.001198093352731686
$ws.Range("G10").Value = 27.77415166666666
$ws.Range("H10").Value = 83.32245499999999
$ws.Range("I10").Value = 0.1759264889641896
$ws.Range("J10").Value = 0.1759264889641896
$ws.Range("O10").Value = 0.02336610738273784
$ws.Range("P10").Value = 0.02336610738273783
$ws.Range("S10").Value = 0.004110717232605298
$ws.Range("T10").Value = 0.004110717232605297
